$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to stay text so values like "1.000" or
# "235.70" are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    'D2' = '30.380.10'
    'E2' = '  +0.09%  '
    'D3' = '1.870.86'
    'E3' = '  -0.70%  '
    'D4' = '1.000'
    'E4' = '  +0.15%  '
    'D5' = '235.70'
    'E5' = '  -0.90%  '
    'E6' = '  +0.20%  '
    'D7' = '0.4668'
    'E7' = '  -0.17%  '
    'D8' = '0.2838'
    'E8' = '  +0.91%  '
    'D9' = '0.06560'
    'E9' = '  -0.04%  '
    'D10' = '21.18'
    'E10' = '  +8.23%  '
    'D11' = '0.07943'
    'E11' = '  +2.73%  '
    'D12' = '97.50'
    'E12' = '  -0.94%  '
    'D13' = '1.866.63'
    'E13' = '  -1.03%  '
    'D14' = '5.154'
    'E14' = '  +0.44%  '
    'E15' = '  +0.88%  '
    'D16' = '282.90'
    'E16' = '  -0.85%  '
    'D17' = '30.381.22'
    'E17' = '  +0.09%  '
    'D18' = '5.555'
    'E18' = '  +4.40%  '
    'E19' = '  +0.18%  '
    'D20' = '12.70'
    'E20' = '  +0.87%  '
    'D21' = '2.114.54'
    'E21' = '  -0.60%  '
    'D22' = '0.000007298'
    'E22' = '  +0.15%  '
    'D23' = '1.001'
    'E23' = '  +0.20%  '
    'D24' = '6.208'
    'E24' = '  +0.28%  '
    'D25' = '9.297'
    'E25' = '  +0.14%  '
    'D26' = '165.03'
    'E26' = '  -1.31%  '
    'D27' = '19.12'
    'E27' = '  +0.48%  '
    'D28' = '1.940'
    'E28' = '  -2.38%  '
    'D29' = '1.355'
    'E29' = '  -1.23%  '
    'D30' = '0.09702'
    'E30' = '  -1.33%  '
    'D31' = '4.439'
    'E32' = '  -1.08%  '
    'D33' = '4.116'
    'E33' = '  -1.72%  '
    'E34' = '  +0.53%  '
    'E35' = '  +2.19%  '
    'D36' = '0.7053'
    'E36' = '  -0.54%  '
    'D37' = '2.719'
    'E37' = '  +0.58%  '
    'D38' = '0.01862'
    'E38' = '  -0.44%  '
    'D39' = '6.327'
    'E39' = '  -4.88%  '
    'D40' = '2.548'
    'E40' = '  +0.90%  '
    'D41' = '73.56'
    'E41' = '  +1.48%  '
    'D42' = '1.948'
    'E42' = '  -1.04%  '
    'B43' = 'TheSandbox'
    'C43' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D43' = '0.4203'
    'E43' = '  +0.13%  '
    'B44' = 'TrustWalletToken'
    'C44' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D44' = '0.8478'
    'E44' = '  -2.58%  '
    'E46' = '  -0.15%  '
    'D47' = '7.219'
    'E47' = '  -0.41%  '
    'D48' = '9.232'
    'E48' = '  -1.82%  '
    'D49' = '940.63'
    'E49' = '  -5.26%  '
    'D50' = '34.21'
    'E50' = '  +0.36%  '
    'E51' = '  -2.45%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# Restore the original (unstyled) look for the cells we touched - the
# text number-format was only needed transiently during assignment.
$ws.Range("D2:E51").Style = "Normal"
